# Feat: txt_debug and reload button on UI
# Feat: SD card detection
# Unfinished: Reload SD card properly
#
# Fills in the two new text-id rows (31 & 32) on the "Translation" sheet
# that were previously just empty placeholder rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Row 31: SingleUseId33 / Default / Right / LTR / <value>
$ws.Range("B31").Value = "SingleUseId33"
$ws.Range("C31").Value = "Default"
$ws.Range("D31").Value = "Right"
$ws.Range("E31").Value = "LTR"
$ws.Range("F31").Value = "<value>"

# Row 32: SingleUseId34 / Default / Left / LTR / <value>
$ws.Range("B32").Value = "SingleUseId34"
$ws.Range("C32").Value = "Default"
$ws.Range("D32").Value = "Left"
$ws.Range("E32").Value = "LTR"
$ws.Range("F32").Value = "<value>"
